# Add new columns I (I0) and J (IF) to Sheet1, mirroring the style of the
# existing header row and populating values for rows 2-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font, thin border, centered) from an existing
# header cell (H1) so the new headers match the rest of the row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data rows (2-29) for columns I and J ---
$data = @{
    2  = @(7, 7)
    3  = @(6, 6)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(8, 8)
    9  = @(9, 9)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(8, 9)
    13 = @(7, 7)
    14 = @(9, 9)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(10, 10)
    23 = @(7, 8)
    24 = @(7, 8)
    25 = @(8, 9)
    26 = @(8, 8)
    27 = @(9, 9)
    28 = @(6, 6)
    29 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # Column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # Column J
}
